$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fixing errors and adding data provider":
# Rename the sheet to reflect its real purpose (invalid login data)
$ws.Name = "invalidLoginDataSheet"

# Remove the extra rows of generated login data (rows 5-11), keeping only
# the first data provider row (rows 1-4: header + 3 data rows)
$ws.Rows("5:11").Delete()

# Leave the selection on D5, as in the saved workbook
$ws.Range("D5").Select()

# Mark the workbook window as minimized
$excel.WindowState = -4140
